# processingPipelineSchema.xlsx update
#
# "zebrafish raw data" gains a new "Friendly Name" column (inserted between
# "Column" and "Description"), and the old placeholder descriptions are
# replaced with real ones. The "Allowable values" column (previously B,
# now shifted to D) keeps its one remaining note.
#
# "zebrafish processed data" gets the same 4-column header row
# (Column / Friendly Name / Description / Allowable values) that the raw
# data sheet uses, replacing the old single stray "xe" cell.

$wb = $excel.ActiveWorkbook

# --- Sheet: "zebrafish raw data" ---
$ws2 = $wb.Worksheets.Item(2)

# Insert a new column B ("Friendly Name"); old B ("Description") becomes C,
# old C ("Allowable values") becomes D.
$ws2.Columns.Item(2).Insert()

$ws2.Range("B1").Value = "Friendly Name"

$ws2.Range("C2").Value = "Internal chemical identifier"
$ws2.Range("C4").Value = "Dose of chemical"
$ws2.Range("C5").Value = "Plate identifier"
$ws2.Range("C6").Value = "Well of plate"
$ws2.Range("C7").Value = "Date assay was completed"
$ws2.Range("C8").Value = "Endpoint measured"
$ws2.Range("C9").Value = "Response (number of fish with phenotype?"

$ws2.Columns.Item(2).EntireColumn.AutoFit()

# --- Sheet: "zebrafish processed data" ---
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A1").Value = "Column"
$ws3.Range("B1").Value = "Friendly Name"
$ws3.Range("C1").Value = "Description"
$ws3.Range("D1").Value = "Allowable values"
$ws3.Range("A1:D1").Font.Bold = $true

# --- Restore/update on-screen selections for each sheet ---
$ws2.Range("C10").Select()
$ws3.Range("A1:D1").Select()

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E3").Select()
